# "front and backend connect"
# Re-sync the room cost-rate table and the raw event log with the latest
# reading reported by the backend.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Room-1 "current rate" column (G) drops by 1.5-2.0 across rows 4-29,
#     and the corresponding room-3 rate cell (O18) ticks up by 0.5.
$ws.Range("G4").Value  = 20
$ws.Range("G5").Value  = 19.5
$ws.Range("G6").Value  = 19
$ws.Range("G7").Value  = 18.5
$ws.Range("G8").Value  = 19
$ws.Range("G9").Value  = 19.5
$ws.Range("G10").Value = 19
$ws.Range("G11").Value = 18.5
$ws.Range("G12").Value = 18
$ws.Range("G13").Value = 18
$ws.Range("G14").Value = 18.5
$ws.Range("G15").Value = 19
$ws.Range("G16").Value = 19.5
$ws.Range("G17").Value = 20
$ws.Range("G18").Value = 20
$ws.Range("G19").Value = 20
$ws.Range("G20").Value = 20
$ws.Range("G21").Value = 20
$ws.Range("G22").Value = 20
$ws.Range("G23").Value = 20
$ws.Range("G24").Value = 20
$ws.Range("G25").Value = 20
$ws.Range("G26").Value = 20
$ws.Range("G27").Value = 20
$ws.Range("G28").Value = 20
$ws.Range("G29").Value = 20

$ws.Range("O18").Value = 24

# --- "费用小计" (cost subtotal) row now has real figures for rooms 1-5.
$ws.Range("B30").Value = "209.667"
$ws.Range("C30").Value = "376.999"
$ws.Range("D30").Value = "155.667"
$ws.Range("E30").Value = "212.667"
$ws.Range("F30").Value = "113.0"

# --- The raw log rows (34-43) were re-captured from the backend a day
#     later; bump every timestamp column (E, I, M, Q, U) to the new run.
$ws.Range("E34").Value = "2024-11-30 21:20:50"
$ws.Range("I34").Value = "2024-11-30 21:20:51"
$ws.Range("M34").Value = "2024-11-30 21:20:52"
$ws.Range("Q34").Value = "2024-11-30 21:20:53"
$ws.Range("U34").Value = "2024-11-30 21:20:51"

$ws.Range("E35").Value = "2024-11-30 21:20:51"
$ws.Range("I35").Value = "2024-11-30 21:20:53"
$ws.Range("M35").Value = "2024-11-30 21:20:53"
$ws.Range("Q35").Value = "2024-11-30 21:20:54"
$ws.Range("U35").Value = "2024-11-30 21:20:54"

$ws.Range("E36").Value = "2024-11-30 21:20:53"
$ws.Range("I36").Value = "2024-11-30 21:20:54"
$ws.Range("M36").Value = "2024-11-30 21:20:56"
$ws.Range("Q36").Value = "2024-11-30 21:20:58"
$ws.Range("U36").Value = "2024-11-30 21:20:55"

$ws.Range("E37").Value = "2024-11-30 21:20:55"
$ws.Range("I37").Value = "2024-11-30 21:20:56"
$ws.Range("M37").Value = "2024-11-30 21:20:58"
$ws.Range("Q37").Value = "2024-11-30 21:20:59"
$ws.Range("U37").Value = "2024-11-30 21:20:56"

$ws.Range("E38").Value = "2024-11-30 21:20:55"
$ws.Range("I38").Value = "2024-11-30 21:20:57"
$ws.Range("M38").Value = "2024-11-30 21:20:58"
$ws.Range("Q38").Value = "2024-11-30 21:20:59"
$ws.Range("U38").Value = "2024-11-30 21:20:57"

$ws.Range("E39").Value = "2024-11-30 21:20:59"
$ws.Range("I39").Value = "2024-11-30 21:20:58"
$ws.Range("M39").Value = "2024-11-30 21:21:02"
$ws.Range("Q39").Value = "2024-11-30 21:20:59"
$ws.Range("U39").Value = "2024-11-30 21:21:02"

$ws.Range("E40").Value = "2024-11-30 21:21:04"
$ws.Range("I40").Value = "2024-11-30 21:20:59"
$ws.Range("M40").Value = "2024-11-30 21:21:02"
$ws.Range("Q40").Value = "2024-11-30 21:21:07"
$ws.Range("U40").Value = "2024-11-30 21:21:02"

$ws.Range("E41").Value = "2024-11-30 21:21:08"
$ws.Range("I41").Value = "2024-11-30 21:21:01"
$ws.Range("M41").Value = "2024-11-30 21:21:03"
$ws.Range("Q41").Value = "2024-11-30 21:21:08"
$ws.Range("U41").Value = "2024-11-30 21:21:03"

$ws.Range("E42").Value = "2024-11-30 21:21:08"
$ws.Range("I42").Value = "2024-11-30 21:21:02"
$ws.Range("M42").Value = "2024-11-30 21:21:04"
$ws.Range("Q42").Value = "2024-11-30 21:21:08"
$ws.Range("U42").Value = "2024-11-30 21:21:03"

$ws.Range("E43").Value = "2024-11-30 21:21:14"
$ws.Range("I43").Value = "2024-11-30 21:21:03"
$ws.Range("M43").Value = "2024-11-30 21:21:04"
$ws.Range("Q43").Value = "2024-11-30 21:21:08"
$ws.Range("U43").Value = "2024-11-30 21:21:05"

# --- Row 44 (room-1 event log entry) was entirely replaced by a fresh
#     "turn off" cycle recorded against the connected backend.
$ws.Range("B44").Value = "turn off"
$ws.Range("D44").Value = 109
$ws.Range("E44").Value = "2024-11-30 21:11:19"
$ws.Range("F44").Value = "turn off"
$ws.Range("G44").Value = 0
$ws.Range("H44").Value = 125.999
$ws.Range("I44").Value = "2024-11-30 21:11:11"
$ws.Range("J44").Value = "add to serve queue"
$ws.Range("L44").Value = 0.99999
$ws.Range("M44").Value = "2024-11-30 21:11:12"
$ws.Range("R44").Value = "turn off"
$ws.Range("S44").Value = 1
$ws.Range("T44").Value = 18.6666
$ws.Range("U44").Value = "2024-11-30 21:11:18"
